$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update score values per the "Authorization Checks - done" evaluation pass
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 31
$ws.Range("C12").Value = 4
$ws.Range("C32").Value = 5

# Move selection / viewport to reflect where the user was working (row 13)
$ws.Range("A7").Select()
$ws.Range("C13").Select()
